$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume strings in column D stay as text (matches original inlineStr cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated values from the diff
$ws.Range("D2").Value = "250.27"
$ws.Range("D4").Value = "5.432"
$ws.Range("D5").Value = "0.05664"
$ws.Range("D6").Value = "3.414"
$ws.Range("D7").Value = "6.379"
$ws.Range("D8").Value = "0.8152"
$ws.Range("D9").Value = "0.9184"
$ws.Range("D10").Value = "0.1439"
$ws.Range("D11").Value = "0.07505"
$ws.Range("D13").Value = "0.03099"
$ws.Range("D14").Value = "0.09353"
$ws.Range("D15").Value = "3.762"
$ws.Range("D16").Value = "0.001589"
$ws.Range("D17").Value = "0.04767"
$ws.Range("D18").Value = "0.0005793"
$ws.Range("D19").Value = "0.006410"
$ws.Range("D20").Value = "0.005036"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D24").Value = "2.178"
$ws.Range("D26").Value = "0.1342"
$ws.Range("D28").Value = "0.0003002"
$ws.Range("D40").Value = "0.04019"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002712"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.006774"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "0.008046"
$ws.Range("D45").Value = "0.00005806"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.5003"
$ws.Range("D49").Value = "0.00002101"
